# Auto-generated edit script: applies per-cell numeric updates (and a few
# cell removals/additions) to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets,
# matching the "Seraph_Profits" price-refresh diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 466.45456
$ws.Range("I2").Value = 297.8
$ws.Range("J2").Value = 607
$ws.Range("K2").Value = 297.8
$ws.Range("L2").Value = 607
$ws.Range("M2").Value = -184.8
$ws.Range("N2").Value = -833

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 2400
$ws.Range("I21").Value = 2400
$ws.Range("K21").Value = 2400
$ws.Range("M21").Value = -1932

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 2400
$ws.Range("I23").Value = 2400
$ws.Range("K23").Value = 2400
$ws.Range("M23").Value = -2166

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 293.58334
$ws.Range("I33").Value = 109.125
$ws.Range("J33").Value = 662.5
$ws.Range("K33").Value = 109.125
$ws.Range("L33").Value = 662.5
$ws.Range("M33").Value = 119.875
$ws.Range("N33").Value = -1120.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 790.875
$ws.Range("I38").Value = 189.57143
$ws.Range("K38").Value = 568.71429
$ws.Range("M38").Value = -196.71429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2221.4443
$ws.Range("J40").Value = 2498.5
$ws.Range("L40").Value = 2498.5
$ws.Range("N40").Value = -2848.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 8500
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 8500
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 8500
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -8638

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2151.3
$ws.Range("J58").Value = 2466.1428
$ws.Range("L58").Value = 7398.428400000001
$ws.Range("N58").Value = -7698.428400000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 8198.25
$ws.Range("J69").Value = 9259.333000000001
$ws.Range("L69").Value = 27777.999
$ws.Range("N69").Value = -29525.999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 8198.25
$ws.Range("J72").Value = 9259.333000000001
$ws.Range("L72").Value = 83333.997
$ws.Range("N72").Value = -92069.997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1011.6667
$ws.Range("I92").Value = 1188.5
$ws.Range("K92").Value = 1188.5
$ws.Range("M92").Value = 59.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2945.75
$ws.Range("I111").Value = 2945.75
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 8837.25
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -5770.25
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6320.7144
$ws.Range("I116").Value = 5840.5713
$ws.Range("J116").Value = 7281
$ws.Range("K116").Value = 5840.5713
$ws.Range("L116").Value = 7281
$ws.Range("M116").Value = -2398.5713
$ws.Range("N116").Value = -14165

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 4539.364
$ws.Range("I131").Value = 1182.0769
$ws.Range("J131").Value = 9388.777
$ws.Range("K131").Value = 3546.2307
$ws.Range("L131").Value = 28166.331
$ws.Range("M131").Value = 1493.7693
$ws.Range("N131").Value = -38246.331

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1912
$ws.Range("I45").Value = 1912
$ws.Range("K45").Value = 1912
$ws.Range("M45").Value = -1535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 241666.67
$ws.Range("J76").Value = 241666.67
$ws.Range("L76").Value = 241666.67
$ws.Range("N76").Value = -242342.67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 241666.67
$ws.Range("J79").Value = 241666.67
$ws.Range("L79").Value = 241666.67
$ws.Range("N79").Value = -244006.67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2672.7144
$ws.Range("I122").Value = 1135.8
$ws.Range("K122").Value = 3407.4
$ws.Range("M122").Value = -957.3999999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1195.7106
$ws.Range("I132").Value = 817.69446
$ws.Range("K132").Value = 2453.08338
$ws.Range("M132").Value = 76.91661999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 499.2353
$ws.Range("I22").Value = 519.1875
$ws.Range("J22").Value = 180
$ws.Range("K22").Value = 519.1875
$ws.Range("L22").Value = 180
$ws.Range("M22").Value = -346.1875
$ws.Range("N22").Value = -526

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2174.8965
$ws.Range("I134").Value = 2074.0356
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 6222.1068
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = -3687.1068
$ws.Range("N134").Value = -20067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2371.5
$ws.Range("I31").Value = 1253.1333
$ws.Range("J31").Value = 4768
$ws.Range("K31").Value = 1253.1333
$ws.Range("L31").Value = 4768
$ws.Range("M31").Value = -958.1333
$ws.Range("N31").Value = -5358

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2371.5
$ws.Range("I34").Value = 1253.1333
$ws.Range("J34").Value = 4768
$ws.Range("K34").Value = 1253.1333
$ws.Range("L34").Value = 4768
$ws.Range("M34").Value = -1051.1333
$ws.Range("N34").Value = -5172

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3792.25
$ws.Range("I105").Value = 1068
$ws.Range("J105").Value = 8332.666999999999
$ws.Range("K105").Value = 1068
$ws.Range("L105").Value = 8332.666999999999
$ws.Range("M105").Value = 679
$ws.Range("N105").Value = -11826.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 6338.0454
$ws.Range("J122").Value = 6162
$ws.Range("L122").Value = 18486
$ws.Range("N122").Value = -23386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 174786
$ws.Range("J141").Value = 174786
$ws.Range("L141").Value = 174786
$ws.Range("N141").Value = -185146

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 135.71428
$ws.Range("I33").Value = 106.666664
$ws.Range("J33").Value = 157.5
$ws.Range("K33").Value = 639.999984
$ws.Range("L33").Value = 945
$ws.Range("M33").Value = -356.999984
$ws.Range("N33").Value = -1511

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1384
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 9000
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 9990.235000000001
$ws.Range("J24").Value = 9990.235000000001
$ws.Range("L24").Value = 9990.235000000001
$ws.Range("N24").Value = -10336.235

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2749.6667
$ws.Range("I80").Value = 1333
$ws.Range("J80").Value = 4166.3335
$ws.Range("K80").Value = 1333
$ws.Range("L80").Value = 4166.3335
$ws.Range("M80").Value = -335
$ws.Range("N80").Value = -6162.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2749.6667
$ws.Range("I83").Value = 1333
$ws.Range("J83").Value = 4166.3335
$ws.Range("K83").Value = 6665
$ws.Range("L83").Value = 20831.6675
$ws.Range("M83").Value = -1673
$ws.Range("N83").Value = -30815.6675

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4401
$ws.Range("I126").Value = 3633
$ws.Range("K126").Value = 10899
$ws.Range("M126").Value = -8429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3824.5
$ws.Range("I16").Value = 3824.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3824.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3654.5
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1014
$ws.Range("I40").Value = 1018
$ws.Range("K40").Value = 1018
$ws.Range("M40").Value = -882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 7498.5
$ws.Range("I100").Value = 997
$ws.Range("J100").Value = 14000
$ws.Range("K100").Value = 997
$ws.Range("L100").Value = 14000
$ws.Range("M100").Value = -456
$ws.Range("N100").Value = -15082

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5579.1577
$ws.Range("I122").Value = 10002
$ws.Range("J122").Value = 5058.8237
$ws.Range("K122").Value = 30006
$ws.Range("L122").Value = 15176.4711
$ws.Range("M122").Value = -27556
$ws.Range("N122").Value = -20076.4711

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 99998.5
$ws.Range("J46").Value = 99998.5
$ws.Range("L46").Value = 99998.5
$ws.Range("N46").Value = -100460.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 99997.664
$ws.Range("J60").Value = 99999.5
$ws.Range("L60").Value = 99999.5
$ws.Range("N60").Value = -101643.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8132.6665
$ws.Range("J62").Value = 8499.357
$ws.Range("L62").Value = 8499.357
$ws.Range("N62").Value = -9747.357

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 65998.336
$ws.Range("J64").Value = 65998.336
$ws.Range("L64").Value = 65998.336
$ws.Range("N64").Value = -66494.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 8132.6665
$ws.Range("J65").Value = 8499.357
$ws.Range("L65").Value = 42496.785
$ws.Range("N65").Value = -48736.785

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 65998.336
$ws.Range("J67").Value = 65998.336
$ws.Range("L67").Value = 65998.336
$ws.Range("N67").Value = -67714.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1349.1666
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1349.1666
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 656.5263
$ws.Range("I122").Value = 592.8125
$ws.Range("K122").Value = 1778.4375
$ws.Range("M122").Value = 671.5625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 99998.5
$ws.Range("J134").Value = 99998.5
$ws.Range("L134").Value = 299995.5
$ws.Range("N134").Value = -305065.5
